$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 190.05882
$ws.Range("I5").Value = 140.46153
$ws.Range("J5").Value = 351.25
$ws.Range("K5").Value = 140.46153
$ws.Range("L5").Value = 351.25
$ws.Range("M5").Value = -25.46153000000001
$ws.Range("N5").Value = -581.25
$ws.Range("H33").Value = 9616868
$ws.Range("I33").Value = 15625547
$ws.Range("J33").Value = 2981.3
$ws.Range("K33").Value = 15625547
$ws.Range("L33").Value = 2981.3
$ws.Range("M33").Value = -15625318
$ws.Range("N33").Value = -3439.3
$ws.Range("H62").Value = 7599.3335
$ws.Range("J62").Value = 7599.3335
$ws.Range("L62").Value = 7599.3335
$ws.Range("N62").Value = -8847.333500000001
$ws.Range("H65").Value = 7599.3335
$ws.Range("J65").Value = 7599.3335
$ws.Range("L65").Value = 37996.6675
$ws.Range("N65").Value = -44236.6675
$ws.Range("H132").Value = 6347.151
$ws.Range("I132").Value = 5231
$ws.Range("K132").Value = 15693
$ws.Range("M132").Value = -13163

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 6313.6665
$ws.Range("I45").Value = 6102.875
$ws.Range("K45").Value = 6102.875
$ws.Range("M45").Value = -5725.875
$ws.Range("H74").Value = 9038.709000000001
$ws.Range("I74").Value = 1653.2094
$ws.Range("K74").Value = 1653.2094
$ws.Range("M74").Value = -779.2094
$ws.Range("H77").Value = 9038.709000000001
$ws.Range("I77").Value = 1653.2094
$ws.Range("K77").Value = 8266.047
$ws.Range("M77").Value = -3898.047
$ws.Range("H102").Value = 18807.584
$ws.Range("I102").Value = 2043.6666
$ws.Range("J102").Value = 69099.336
$ws.Range("K102").Value = 2043.6666
$ws.Range("L102").Value = 69099.336
$ws.Range("M102").Value = -421.6666
$ws.Range("N102").Value = -72343.336
$ws.Range("H133").Value = 57259.2
$ws.Range("J133").Value = 57259.2
$ws.Range("L133").Value = 57259.2
$ws.Range("N133").Value = -62319.2
$ws.Range("H134").Value = 42985.57
$ws.Range("J134").Value = 42985.57
$ws.Range("L134").Value = 42985.57
$ws.Range("N134").Value = -53125.57

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 11523471
$ws.Range("I20").Value = 25650586
$ws.Range("J20").Value = 45188.938
$ws.Range("K20").Value = 25650586
$ws.Range("L20").Value = 45188.938
$ws.Range("M20").Value = -25650339
$ws.Range("N20").Value = -45682.938
$ws.Range("H22").Value = 11910372
$ws.Range("I22").Value = 12992678
$ws.Range("J22").Value = 5000
$ws.Range("K22").Value = 12992678
$ws.Range("L22").Value = 5000
$ws.Range("M22").Value = -12992505
$ws.Range("N22").Value = -5346
$ws.Range("H96").Value = 34500
$ws.Range("I96").Value = 22100
$ws.Range("J96").Value = 50000
$ws.Range("K96").Value = 22100
$ws.Range("L96").Value = 50000
$ws.Range("M96").Value = -19354
$ws.Range("N96").Value = -55492
$ws.Range("H107").Value = 3426.125
$ws.Range("I107").Value = 3901.5
$ws.Range("J107").Value = 2000
$ws.Range("K107").Value = 3901.5
$ws.Range("L107").Value = 2000
$ws.Range("M107").Value = -1981.5
$ws.Range("N107").Value = -5840

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 125344.75
$ws.Range("J10").Value = 999999
$ws.Range("L10").Value = 999999
$ws.Range("N10").Value = -1000277
$ws.Range("H31").Value = 70647.086
$ws.Range("I31").Value = 141619.94
$ws.Range("J31").Value = 17417.45
$ws.Range("K31").Value = 141619.94
$ws.Range("L31").Value = 17417.45
$ws.Range("M31").Value = -141324.94
$ws.Range("N31").Value = -18007.45
$ws.Range("H34").Value = 70647.086
$ws.Range("I34").Value = 141619.94
$ws.Range("J34").Value = 17417.45
$ws.Range("K34").Value = 141619.94
$ws.Range("L34").Value = 17417.45
$ws.Range("M34").Value = -141417.94
$ws.Range("N34").Value = -17821.45
$ws.Range("H58").Value = 13598.667
$ws.Range("I58").Value = 4388.1577
$ws.Range("K58").Value = 4388.1577
$ws.Range("M58").Value = -4185.1577
$ws.Range("H62").Value = 4313.5713
$ws.Range("I62").Value = 4765
$ws.Range("K62").Value = 4765
$ws.Range("M62").Value = -4141
$ws.Range("H65").Value = 4313.5713
$ws.Range("I65").Value = 4765
$ws.Range("K65").Value = 23825
$ws.Range("M65").Value = -20705
$ws.Range("H70").Value = 17000
$ws.Range("J70").Value = 17000
$ws.Range("L70").Value = 17000
$ws.Range("N70").Value = -17630
$ws.Range("H73").Value = 17000
$ws.Range("J73").Value = 17000
$ws.Range("L73").Value = 17000
$ws.Range("N73").Value = -19184
$ws.Range("H136").Value = 13598.667
$ws.Range("I136").Value = 4388.1577
$ws.Range("K136").Value = 13164.4731
$ws.Range("M136").Value = -10614.4731

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 876.76666
$ws.Range("J107").Value = 1102.5333
$ws.Range("L107").Value = 3307.5999
$ws.Range("N107").Value = -7147.5999
$ws.Range("H136").Value = 556.1429000000001
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()
$ws.Range("H137").Value = 3430
$ws.Range("I137").Value = 2150
$ws.Range("J137").Value = 5990
$ws.Range("K137").Value = 6450
$ws.Range("L137").Value = 17970
$ws.Range("M137").Value = -1350
$ws.Range("N137").Value = -28170
$ws.Range("H140").Value = 2003.5714
$ws.Range("I140").Value = 2003.5714
$ws.Range("J140").Value = 0
$ws.Range("K140").Value = 6010.7142
$ws.Range("L140").Value = 0
$ws.Range("N140").Value = -830.7142000000003

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 23262
$ws.Range("I80").Value = 31444
$ws.Range("J80").Value = 20924.285
$ws.Range("K80").Value = 31444
$ws.Range("L80").Value = 20924.285
$ws.Range("M80").Value = -30446
$ws.Range("N80").Value = -22920.285
$ws.Range("H83").Value = 23262
$ws.Range("I83").Value = 31444
$ws.Range("J83").Value = 20924.285
$ws.Range("K83").Value = 157220
$ws.Range("L83").Value = 104621.425
$ws.Range("M83").Value = -152228
$ws.Range("N83").Value = -114605.425
$ws.Range("H97").Value = 3195.5264
$ws.Range("I97").Value = 1426.6666
$ws.Range("K97").Value = 1426.6666
$ws.Range("M97").Value = -930.6666
$ws.Range("H105").Value = 79118.625
$ws.Range("J105").Value = 79118.625
$ws.Range("L105").Value = 79118.625
$ws.Range("N105").Value = -86106.625
$ws.Range("H126").Value = 4885156
$ws.Range("I126").Value = 3596022.2
$ws.Range("J126").Value = 5959434
$ws.Range("K126").Value = 10788066.6
$ws.Range("L126").Value = 17878302
$ws.Range("M126").Value = -10785596.6
$ws.Range("N126").Value = -17883242

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 100002424
$ws.Range("I16").Value = 100002424
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 100002424
$ws.Range("L16").Value = 0
$ws.Range("N16").Value = -100002254
$ws.Range("H100").Value = 2775.647
$ws.Range("I100").Value = 2616.5833
$ws.Range("K100").Value = 2616.5833
$ws.Range("M100").Value = -2075.5833
$ws.Range("H132").Value = 6598.4717
$ws.Range("J132").Value = 23675.6
$ws.Range("L132").Value = 71026.79999999999
$ws.Range("N132").Value = -76086.79999999999
$ws.Range("H136").Value = 47346.043
$ws.Range("I136").Value = 49657.383
$ws.Range("K136").Value = 148972.149
$ws.Range("M136").Value = -146422.149

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 132165
$ws.Range("I62").Value = 21720
$ws.Range("J62").Value = 242610
$ws.Range("K62").Value = 21720
$ws.Range("L62").Value = 242610
$ws.Range("M62").Value = -21096
$ws.Range("N62").Value = -243858
$ws.Range("H65").Value = 132165
$ws.Range("I65").Value = 21720
$ws.Range("J65").Value = 242610
$ws.Range("K65").Value = 108600
$ws.Range("L65").Value = 1213050
$ws.Range("M65").Value = -105480
$ws.Range("N65").Value = -1219290
$ws.Range("H132").Value = 6055.2764
$ws.Range("I132").Value = 2170.926
$ws.Range("K132").Value = 6512.778
$ws.Range("M132").Value = -3982.778
$ws.Range("H136").Value = 13231.429
$ws.Range("I136").Value = 1353.4615
$ws.Range("K136").Value = 4060.3845
$ws.Range("M136").Value = -1510.3845
